# Update odds values for row 27 on the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G27"  = 1.53
    "H27"  = 4.2
    "I27"  = 5.5
    "J27"  = 2.05
    "K27"  = 2.4
    "L27"  = 5.5
    "U27"  = 1.3
    "V27"  = 3.4
    "W27"  = 1.75
    "X27"  = 2
    "Y27"  = 8
    "Z27"  = 8
    "AB27" = 11
    "AD27" = 23
    "AE27" = 15
    "AG27" = 17
    "AI27" = 201
    "AK27" = 29
    "AL27" = 17
    "AM27" = 51
    "AN27" = 41
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
